$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = "66.695.70"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.785.58"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "432.91"
$ws.Range("E5").Value = "  +5.62%  "
$ws.Range("D6").Value = "140.13"
$ws.Range("E6").Value = "  +7.23%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.734"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  -9.60%  "
$ws.Range("D11").Value = "0.0000313"
$ws.Range("E11").Value = "  -13.78%  "
$ws.Range("D12").Value = "42.82"
$ws.Range("E12").Value = "  +4.83%  "
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "4.385.78"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "14.87"
$ws.Range("E15").Value = "  -5.93%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "3.789.62"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("E19").Value = "  +7.14%  "
$ws.Range("D20").Value = "66.740.04"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "409.40"
$ws.Range("D22").Value = "14.75"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").Value = "3.27"
$ws.Range("E23").Value = "  +7.32%  "
$ws.Range("D24").Value = "85.59"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "36.86"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +7.80%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "9.77"
$ws.Range("E28").Value = "  +36.50%  "
$ws.Range("D29").Value = "9.78"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "13.87"
$ws.Range("E30").Value = "  +11.20%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "712.84"
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("E32").Value = "  +10.28%  "
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "42.03"
$ws.Range("E34").Value = "  +8.98%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +26.90%  "
$ws.Range("D38").Value = "56.07"
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  +40.89%  "
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("D42").Value = "3.37"
$ws.Range("E42").Value = "  +8.27%  "
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").Value = "0.0₃0675"
$ws.Range("E44").Value = "  -14.08%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "0.324"
$ws.Range("E46").Value = "  +9.85%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").Value = "2.08"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").Value = "142.50"
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("E51").Value = "  +1.72%  "
